$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 3007.1785
$ws.Range("J43").Value = 3083.3333
$ws.Range("L43").Value = 3083.3333
$ws.Range("N43").Value = -3221.3333

$ws.Range("H98").Value = 800.3077
$ws.Range("I98").Value = 783.6667
$ws.Range("J98").Value = 1000
$ws.Range("K98").Value = 783.6667
$ws.Range("L98").Value = 1000
$ws.Range("M98").Value = 714.3333
$ws.Range("N98").Value = -3996

$ws.Range("H103").Value = 2572.3333
$ws.Range("J103").Value = 5150.5
$ws.Range("L103").Value = 15451.5
$ws.Range("N103").Value = -16623.5

$ws.Range("H108").Value = 39787.668
$ws.Range("J108").Value = 39787.668
$ws.Range("L108").Value = 39787.668
$ws.Range("N108").Value = -47467.668

$ws.Range("H113").Value = 501252.5
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()

$ws.Range("H116").Value = 2152.875
$ws.Range("J116").Value = 1875.25
$ws.Range("L116").Value = 1875.25
$ws.Range("N116").Value = -8759.25

$ws.Range("H122").Value = 800.3077
$ws.Range("I122").Value = 783.6667
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 2351.0001
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = 98.9998999999998
$ws.Range("N122").Value = -7900

$ws.Range("H129").Value = 2654.6843
$ws.Range("I129").Value = 14914.286
$ws.Range("J129").Value = 938.34
$ws.Range("K129").Value = 44742.858
$ws.Range("L129").Value = 2815.02
$ws.Range("M129").Value = -39742.858
$ws.Range("N129").Value = -12815.02

$ws.Range("H132").Value = 4242481.5
$ws.Range("I132").Value = 5005414.5
$ws.Range("J132").Value = 3967.111
$ws.Range("K132").Value = 15016243.5
$ws.Range("L132").Value = 11901.333
$ws.Range("M132").Value = -15013713.5
$ws.Range("N132").Value = -16961.333

$ws.Range("H137").Value = 1486.1666
$ws.Range("I137").Value = 1514.7646
$ws.Range("J137").Value = 1000
$ws.Range("K137").Value = 4544.293799999999
$ws.Range("L137").Value = 3000
$ws.Range("M137").Value = -1994.293799999999
$ws.Range("N137").Value = -8100

$ws.Range("H138").Value = 3433.8276
$ws.Range("I138").Value = 3871.2222
$ws.Range("J138").Value = 3383.359
$ws.Range("K138").Value = 11613.6666
$ws.Range("L138").Value = 10150.077
$ws.Range("M138").Value = -6473.6666
$ws.Range("N138").Value = -20430.077

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 29154.314
$ws.Range("I32").Value = 13220.779
$ws.Range("J32").Value = 60490.266
$ws.Range("K32").Value = 13220.779
$ws.Range("L32").Value = 60490.266
$ws.Range("M32").Value = -12933.779
$ws.Range("N32").Value = -61064.266

$ws.Range("H74").Value = 1092.8334
$ws.Range("I74").Value = 674.2857
$ws.Range("J74").Value = 1678.8
$ws.Range("K74").Value = 674.2857
$ws.Range("L74").Value = 1678.8
$ws.Range("M74").Value = 199.7143
$ws.Range("N74").Value = -3426.8

$ws.Range("H77").Value = 1092.8334
$ws.Range("I77").Value = 674.2857
$ws.Range("J77").Value = 1678.8
$ws.Range("K77").Value = 3371.4285
$ws.Range("L77").Value = 8394
$ws.Range("M77").Value = 996.5715
$ws.Range("N77").Value = -17130

$ws.Range("H122").Value = 3097
$ws.Range("I122").Value = 2711.25
$ws.Range("K122").Value = 8133.75
$ws.Range("M122").Value = -5683.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 225801.2
$ws.Range("I86").Value = 225801.2
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 225801.2
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -224678.2
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 225801.2
$ws.Range("I89").Value = 225801.2
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 1129006
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -1123390
$ws.Range("N89").ClearContents()

$ws.Range("H134").Value = 3987.653
$ws.Range("I134").Value = 4424.5713
$ws.Range("J134").Value = 2895.3572
$ws.Range("K134").Value = 13273.7139
$ws.Range("L134").Value = 8686.071599999999
$ws.Range("M134").Value = -10738.7139
$ws.Range("N134").Value = -13756.0716

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1743.0588
$ws.Range("I134").Value = 986.75
$ws.Range("J134").Value = 2415.3333
$ws.Range("K134").Value = 2960.25
$ws.Range("L134").Value = 7245.999899999999
$ws.Range("M134").Value = -425.25
$ws.Range("N134").Value = -12315.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1626.2693
$ws.Range("I5").Value = 1398.7858
$ws.Range("J5").Value = 1891.6666
$ws.Range("K5").Value = 4196.357400000001
$ws.Range("L5").Value = 5674.9998
$ws.Range("M5").Value = -4084.357400000001
$ws.Range("N5").Value = -5898.9998

$ws.Range("H19").Value = 1000
$ws.Range("J19").Value = 1000
$ws.Range("L19").Value = 3000
$ws.Range("N19").Value = -3348

$ws.Range("H22").Value = 3094.0981
$ws.Range("J22").Value = 3179.5715
$ws.Range("L22").Value = 9538.7145
$ws.Range("N22").Value = -9876.7145

$ws.Range("H27").Value = 3094.0981
$ws.Range("J27").Value = 3179.5715
$ws.Range("L27").Value = 9538.7145
$ws.Range("N27").Value = -9742.7145

$ws.Range("H131").Value = 642805.7
$ws.Range("I131").Value = 700
$ws.Range("J131").Value = 726558.5600000001
$ws.Range("K131").Value = 2100
$ws.Range("L131").Value = 2179675.68
$ws.Range("M131").Value = 2940
$ws.Range("N131").Value = -2189755.68

$ws.Range("H135").Value = 1626.2693
$ws.Range("I135").Value = 1398.7858
$ws.Range("J135").Value = 1891.6666
$ws.Range("K135").Value = 12589.0722
$ws.Range("L135").Value = 17024.9994
$ws.Range("M135").Value = -10054.0722
$ws.Range("N135").Value = -22094.9994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 44800
$ws.Range("J20").Value = 44800
$ws.Range("L20").Value = 44800
$ws.Range("N20").Value = -45290

$ws.Range("H24").Value = 8000
$ws.Range("J24").Value = 8000
$ws.Range("L24").Value = 8000
$ws.Range("N24").Value = -8346

$ws.Range("H43").Value = 1975.0646
$ws.Range("J43").Value = 7003.8
$ws.Range("L43").Value = 7003.8
$ws.Range("N43").Value = -7305.8

$ws.Range("H46").Value = 12486.6
$ws.Range("J46").Value = 12486.6
$ws.Range("L46").Value = 12486.6
$ws.Range("N46").Value = -12798.6

$ws.Range("H57").Value = 19500
$ws.Range("J57").Value = 19500
$ws.Range("L57").Value = 19500
$ws.Range("N57").Value = -21140

$ws.Range("H80").Value = 143173280
$ws.Range("I80").Value = 200441600
$ws.Range("J80").Value = 2500
$ws.Range("K80").Value = 200441600
$ws.Range("L80").Value = 2500
$ws.Range("M80").Value = -200440602
$ws.Range("N80").Value = -4496

$ws.Range("H83").Value = 143173280
$ws.Range("I83").Value = 200441600
$ws.Range("J83").Value = 2500
$ws.Range("K83").Value = 1002208000
$ws.Range("L83").Value = 12500
$ws.Range("M83").Value = -1002203008
$ws.Range("N83").Value = -22484

$ws.Range("H126").Value = 3302.4707
$ws.Range("I126").Value = 3248.4167
$ws.Range("J126").Value = 3432.2
$ws.Range("K126").Value = 9745.250100000001
$ws.Range("L126").Value = 10296.6
$ws.Range("M126").Value = -7275.250100000001
$ws.Range("N126").Value = -15236.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2715.3333
$ws.Range("I7").Value = 1676.8
$ws.Range("J7").Value = 3457.1428
$ws.Range("K7").Value = 1676.8
$ws.Range("L7").Value = 3457.1428
$ws.Range("M7").Value = -1564.8
$ws.Range("N7").Value = -3681.1428

$ws.Range("H55").Value = 1485.9395
$ws.Range("J55").Value = 1400.2609
$ws.Range("L55").Value = 1400.2609
$ws.Range("N55").Value = -1746.2609

$ws.Range("H126").Value = 2715.3333
$ws.Range("I126").Value = 1676.8
$ws.Range("J126").Value = 3457.1428
$ws.Range("K126").Value = 5030.4
$ws.Range("L126").Value = 10371.4284
$ws.Range("M126").Value = -2560.4
$ws.Range("N126").Value = -15311.4284

$ws.Range("H132").Value = 4270.926
$ws.Range("I132").Value = 6098.9165
$ws.Range("K132").Value = 18296.7495
$ws.Range("M132").Value = -15766.7495

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 6818.647
$ws.Range("I54").Value = 6380
$ws.Range("J54").Value = 6912.643
$ws.Range("K54").Value = 6380
$ws.Range("L54").Value = 6912.643
$ws.Range("M54").Value = -5860
$ws.Range("N54").Value = -7952.643

$ws.Range("H81").Value = 183505.36
$ws.Range("J81").Value = 144922.58
$ws.Range("L81").Value = 289845.16
$ws.Range("N81").Value = -291967.16

$ws.Range("H84").Value = 183505.36
$ws.Range("J84").Value = 144922.58
$ws.Range("L84").Value = 1449225.8
$ws.Range("N84").Value = -1459833.8

$ws.Range("H107").Value = 167162.83
$ws.Range("I107").Value = 526.6667
$ws.Range("J107").Value = 333799
$ws.Range("K107").Value = 1580.0001
$ws.Range("L107").Value = 1001397
$ws.Range("M107").Value = 339.9999
$ws.Range("N107").Value = -1005237

$ws.Range("H126").Value = 1269.8182
$ws.Range("I126").Value = 1372.25
$ws.Range("J126").Value = 996.6667
$ws.Range("K126").Value = 4116.75
$ws.Range("L126").Value = 2990.0001
$ws.Range("M126").Value = -1646.75
$ws.Range("N126").Value = -7930.0001

$ws.Range("H136").Value = 3286.5
$ws.Range("I136").Value = 4039.1035
$ws.Range("J136").Value = 2478.1482
$ws.Range("K136").Value = 12117.3105
$ws.Range("L136").Value = 7434.444600000001
$ws.Range("M136").Value = -9567.3105
$ws.Range("N136").Value = -12534.4446
